$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1/G1 hold numeric-looking rating values that must be stored as text (like the source data),
# so force a text number format before assigning them.
$ws.Range("F1:G1").NumberFormat = "@"

$ws.Range("A1").Value = "Caterpillar Colorado Erkek Kahverengi Nubuk Bot-637"
$ws.Range("D1").Value = "https://n11scdn.akamaized.net/a1/602_857/13/22/22/80/IMG-4041382771724777596.jpg"
$ws.Range("E1").Value = "['https://n11scdn.akamaized.net/a1/602_857/13/22/22/80/IMG-4041382771724777596.jpg', 'https://n11scdn.akamaized.net/a1/60_86/13/22/22/80/IMG-4041382771724777596.jpg', 'https://n11scdn.akamaized.net/a1/60_86/08/64/63/42/IMG-6424608762558479759.jpg', 'https://n11scdn.akamaized.net/a1/60_86/07/88/62/46/IMG-8736130210984762394.jpg', 'https://n11scdn.akamaized.net/a1/60_86/01/83/99/49/IMG-489356320240617610.jpg', 'https://n11scdn.akamaized.net/a1/60_86/09/61/74/63/IMG-9115817892814919145.jpg', 'https://n11scdn.akamaized.net/a1/500_700/24/01/17/72/20/82/70/41/44/34/40/23/97881136301522895152.png']"
$ws.Range("F1").Value = "4.5"
$ws.Range("G1").Value = "5"
